$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price column so numeric-looking values
# (e.g. "1.003") are stored as text, matching the source data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.205.87"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.660.30"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "217.30"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "0.5191"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.2640"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "0.06277"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "20.82"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "0.07780"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "4.477"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.664.04"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "1.885.87"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "0.5479"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "0.0₅8118"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "64.95"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "26.206.67"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "4.627"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "192.44"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "6.010"
$ws.Range("E23").Value = "  -4.51%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "139.62"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "0.1222"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").Value = "7.300"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "1.436"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "0.05960"
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("D31").Value = "1.273"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "3.552"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("D33").Value = "3.281"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").Value = "0.9623"
$ws.Range("E35").Value = "  -4.13%  "
$ws.Range("D36").Value = "2.419"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "2.769"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "0.5699"
$ws.Range("E38").Value = "  -6.03%  "
$ws.Range("D39").Value = "6.026"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "0.8493"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "1.005.14"
$ws.Range("E43").Value = "  -7.89%  "
$ws.Range("D44").Value = "100.37"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "1.800.65"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "56.58"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "8.019"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "0.4315"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "0.05167"
$ws.Range("E51").Value = "  -0.77%  "

# Restore the original (default) style now that values are set as text.
$priceRange.Style = "Normal"
